$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "supplier_id" column (A) is dropped from the export: the remaining
# columns (supplier_kode / supplier_nama / supplier_alamat) each shift one
# column to the left, and the old column D is emptied out.

# Header row.
$ws.Range("A1").Value = "supplier_kode"
$ws.Range("B1").Value = "supplier_nama"
$ws.Range("C1").Value = "supplier_alamat"

# Row 2 (was SUP004 / UD. Sumber Rejeki / Jl. Bunga No. 20, Surabaya).
$ws.Range("A2").Value = "SUP004"
$ws.Range("B2").Value = "UD. Sumber Rejeki"
$ws.Range("C2").Value = "Jl. Bunga No. 20, Surabaya"

# Row 3 (was SUP005 / PT. Cipta Karya / Jl. Cempaka No. 15, Yogyakarta).
$ws.Range("A3").Value = "SUP005"
$ws.Range("B3").Value = "PT. Cipta Karya"
$ws.Range("C3").Value = "Jl. Cempaka No. 15, Yogyakarta"

# Columns B/C in rows 2-3 now hold what used to live in C/D, which used the
# plain bordered style (no wrap / vertical-center) instead of the wrapped
# style the B column used to carry - copy that formatting over.
$ws.Range("C2").Copy()
$ws.Range("B2:B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The old column D is no longer used.
$ws.Range("D1:D3").Clear()

# Match the new column widths recorded after the edit.
$ws.Range("B1").ColumnWidth = 15.3
$ws.Range("C1").ColumnWidth = 29.5

# Selection recorded after the edit.
$ws.Range("D3").Select()

$wb.Save()
